$wb = $excel.ActiveWorkbook

# --- Sheet "Resource Utilization": update the two utilization figures ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = 99.77
$ws1.Range("B3").Value = 39.55

# --- Sheet "Activity Times": the simulation report is now missing most
#     of its detail, so the per-step breakdown rows (3-8) disappear and
#     only the summary row (row 2) remains, with most of its metrics blank.
$ws2 = $wb.Worksheets.Item(2)

# Drop the detailed activity rows entirely.
$ws2.Range("A3:K8").EntireRow.Delete()

# Refresh the summary row: tokens counts collapse to 0, and the timing
# columns are no longer available.
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0

# Leave E2:G2 present but blank (matching the pre-existing blank cells
# H2:K2) rather than deleting them outright. A bare "'" is Excel's
# force-text prefix for an empty string, so the cell keeps real (blank)
# content instead of becoming a genuinely empty/absent cell; resetting
# the style afterwards drops the cosmetic quote-prefix formatting it
# would otherwise leave behind.
$ws2.Range("E2").Value = "'"
$ws2.Range("E2").Style = "Normal"
$ws2.Range("F2").Value = "'"
$ws2.Range("F2").Style = "Normal"
$ws2.Range("G2").Value = "'"
$ws2.Range("G2").Style = "Normal"
